$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.446.17"
$ws.Range("E2").Value = "  +1.60%  "
$ws.Range("D3").Value = "3.775.86"
$ws.Range("E3").Value = "  +0.64%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'596.40"
$ws.Range("E5").Value = "  +0.45%  "
$ws.Range("D6").Value = "'168.74"
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").Value = "3.775.82"
$ws.Range("E7").Value = "  +0.59%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("E9").Value = "  -0.42%  "
$ws.Range("E10").Value = "  -0.88%  "
$ws.Range("D11").Value = "'6.53"
$ws.Range("E11").Value = "  +1.02%  "
$ws.Range("D12").Value = "'0.451"
$ws.Range("E12").Value = "  -0.53%  "
$ws.Range("D13").Value = "'0.0000265"
$ws.Range("E13").Value = "  -2.43%  "
$ws.Range("D14").Value = "'36.81"
$ws.Range("E14").Value = "  +0.82%  "
$ws.Range("D15").Value = "4.407.11"
$ws.Range("E15").Value = "  +0.27%  "
$ws.Range("D16").Value = "3.780.06"
$ws.Range("E16").Value = "  +0.55%  "
$ws.Range("D17").Value = "68.442.29"
$ws.Range("E17").Value = "  +1.38%  "
$ws.Range("D18").Value = "'18.22"
$ws.Range("E18").Value = "  -3.18%  "
$ws.Range("D19").Value = "'7.05"
$ws.Range("E19").Value = "  -1.88%  "
$ws.Range("E20").Value = "  -0.09%  "
$ws.Range("E21").Value = "  +4.10%  "
$ws.Range("D22").Value = "'468.39"
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("E23").Value = "  -2.59%  "
$ws.Range("D24").Value = "'85.18"
$ws.Range("E24").Value = "  +1.82%  "
$ws.Range("D25").Value = "'0.0000145"
$ws.Range("E25").Value = "  -1.51%  "
$ws.Range("E26").Value = "  +0.95%  "
$ws.Range("D27").Value = "'12.23"
$ws.Range("E27").Value = "  +0.78%  "
$ws.Range("D28").Value = "'10.18"
$ws.Range("E28").Value = "  -0.52%  "
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("D30").Value = "3.920.53"
$ws.Range("E30").Value = "  +0.44%  "
$ws.Range("E31").Value = "  -2.92%  "
$ws.Range("D32").Value = "'7.40"
$ws.Range("E32").Value = "  -2.96%  "
$ws.Range("E33").Value = "  -0.82%  "
$ws.Range("D34").Value = "'30.12"
$ws.Range("E34").Value = "  -0.72%  "
$ws.Range("D35").Value = "'9.36"
$ws.Range("E35").Value = "  +2.50%  "
$ws.Range("D36").Value = "'0.998"
$ws.Range("D37").Value = "3.727.13"
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("E38").Value = "  -2.47%  "
$ws.Range("D39").Value = "'3.49"
$ws.Range("E39").Value = "  -9.00%  "
$ws.Range("E40").Value = "  +1.31%  "
$ws.Range("E41").Value = "  +1.06%  "
$ws.Range("D42").Value = "'5.83"
$ws.Range("E42").Value = "  -0.67%  "
$ws.Range("D43").Value = "'0.999"
$ws.Range("E43").Value = "  -0.18%  "
$ws.Range("D44").Value = "'0.310"
$ws.Range("E44").Value = "  -0.94%  "
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").Value = "'1.96"
$ws.Range("E46").Value = "  +0.83%  "
$ws.Range("D47").Value = "'8.60"
$ws.Range("E47").Value = "  -0.96%  "
$ws.Range("D48").Value = "'42.53"
$ws.Range("E48").Value = "  +10.16%  "
$ws.Range("D49").Value = "'403.35"
$ws.Range("E49").Value = "  +1.84%  "
$ws.Range("D50").Value = "'45.62"
$ws.Range("E50").Value = "  -0.16%  "
$ws.Range("D51").Value = "'145.96"
$ws.Range("E51").Value = "  +3.13%  "
